# Apply updated Betfair Back/Lay odds values to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3" = 2.14
    "W3" = 1.73
    "I6" = 2.44
    "J6" = 3.3
    "N6" = 3
    "P6" = 1.68
    "Q6" = 2.22
    "V6" = 1.69
    "F7" = 1.65
    "J7" = 3.6
    "K7" = 5.2
    "H9" = 2.36
    "K9" = 4
    "L9" = 1.29
    "P9" = 2.08
    "F10" = 1.79
    "X10" = 15
    "H12" = 2.34
    "M12" = 1.12
    "N12" = 2.98
    "P12" = 1.65
    "Q12" = 2.48
    "T12" = 2.08
    "U12" = 1.88
    "X12" = 9.4
    "AO12" = 28
    "Y13" = 16
    "AC13" = 9
    "AG13" = 11.5
    "AN13" = 14.5
    "F14" = 2.36
    "H14" = 3.05
    "I14" = 3.45
    "Z14" = 980
    "AA14" = 60
    "AO14" = 980
    "L15" = 1.42
    "N15" = 3.45
    "G16" = 2.76
    "H16" = 2.6
    "J16" = 3.8
    "K16" = 4.3
    "L16" = 1.31
    "P16" = 2.34
    "R16" = 1.54
    "U16" = 2.48
    "W16" = 1.57
    "F17" = 1.95
    "G17" = 2.1
    "H17" = 3.6
    "J17" = 3.85
    "N17" = 4.9
    "P17" = 2.34
    "Q17" = 1.63
    "R17" = 1.53
    "S17" = 2.56
    "W17" = 1.9
    "AG17" = 13
    "AJ17" = 29
    "AK17" = 23
    "AN17" = 12.5
    "F18" = 1.62
    "AE18" = 60
    "AM18" = 75
    "AN18" = 7
    "H19" = 3.9
    "L19" = 1.29
    "N19" = 5.1
    "U19" = 2.36
    "G20" = 3
    "Q20" = 2.58
    "R20" = 1.2
    "W20" = 1.5
    "AB20" = 9.800000000000001
    "AC20" = 8
    "AG20" = 14
    "F21" = 2.48
    "N21" = 2.78
    "T21" = 1.98
    "X21" = 10.5
    "G22" = 4.4
    "I22" = 2.2
    "U22" = 1.92
    "V22" = 1.83
    "AN22" = 80
    "F23" = 2.14
    "G23" = 2.16
    "W23" = 1.86
    "Z23" = 28
    "L24" = 1.69
    "F25" = 1.8
    "G25" = 1.81
    "H25" = 5
    "I25" = 5.1
    "R25" = 1.49
    "S25" = 3
    "T25" = 1.76
    "V25" = 1.24
    "W25" = 2.22
    "AD25" = 19.5
    "AH25" = 18.5
    "AK25" = 18.5
    "AO25" = 60
    "F26" = 2.12
    "G26" = 2.14
    "I26" = 4.3
    "N26" = 3.25
    "T26" = 1.97
    "V26" = 1.3
    "W26" = 1.88
    "Y26" = 13
    "AB26" = 8
    "AJ26" = 25
    "P27" = 1.75
    "AH27" = 26
    "AN27" = 13.5
    "J28" = 1.03
    "F29" = 2.84
    "Q29" = 2.34
    "F30" = 2.12
    "G30" = 2.34
    "I30" = 4.1
    "Q30" = 1.96
    "V30" = 1.33
    "H31" = 3.6
    "I31" = 3.95
    "M31" = 1.09
    "Q31" = 2.14
    "V31" = 1.34
    "K32" = 3.7
    "L32" = 1.56
    "Q32" = 2.6
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
